$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: titles (shared strings with preserved whitespace) ---
$ws.Range("A2").Value = "`n           Feriantes`n          "
$ws.Range("A3").Value = "`n           Everywoman`n          "
$ws.Range("A4").Value = "`n           Misericordia`n          "
$ws.Range("A5").Value = "`n           Ricardo III`n          "
$ws.Range("A6").Value = "`n             Sans tambour`n            "
$ws.Range("A7").Value = "`n           IRIBARNE`n          "
$ws.Range("A8").Value = "`n           Obra infinita`n          "
$ws.Range("A9").Value = "`n           Pieces of a Woman`n          "
$ws.Range("A10").Value = "`n            Hamlet`n           "
$ws.Range("A11").Value = "`n           Dragón`n          "
$ws.Range("A12").Value = "`n           Comedia sin título`n          "
$ws.Range("A13").Value = "`n           Imitation of Life`n          "
$ws.Range("A14").Value = "`n           El patito feo`n          "
$ws.Range("A15").Value = "`n           Colgando de un hilo`n          "

# --- Column B: dates, using mm-dd-yy which maps to built-in numFmtId 14 ---
$ws.Range("B2").Value = 45276
$ws.Range("B2").NumberFormat = "mm-dd-yy"
$ws.Range("B2").Copy($ws.Range("B3:B15"))
$ws.Range("B3").Value = 45304
$ws.Range("B4").Value = 45336
$ws.Range("B5").Value = 45255
$ws.Range("B6").Value = 45247
$ws.Range("B7").Value = 45212
$ws.Range("B8").Value = 44992
$ws.Range("B9").Value = 44910
$ws.Range("B10").Value = 44729
$ws.Range("B11").Value = 44605
$ws.Range("B12").Value = 44518
$ws.Range("B13").Value = 44521
$ws.Range("B14").Value = 44471
$ws.Range("B15").Value = 45258

# --- Column C: ratings ---
$ws.Range("C2").Value = 7
$ws.Range("C3").Value = 5
$ws.Range("C4").Value = 10
$ws.Range("C5").Value = 6
$ws.Range("C6").Value = 3
$ws.Range("C7").Value = 9
$ws.Range("C8").Value = 4
$ws.Range("C9").Value = 8
$ws.Range("C10").Value = 10
$ws.Range("C11").Value = 7
$ws.Range("C12").Value = 8
$ws.Range("C13").Value = 2
$ws.Range("C14").Value = 5
$ws.Range("C15").Value = 7

# --- Column width for A, matching the recorded original as closely as possible ---
$ws.Columns.Item(1).ColumnWidth = 36.5

# --- Selection matches final author session state ---
$ws.Range("C16").Select() | Out-Null
